# Week 5 ToDo.docx edit script
# Commit message: "Made Progress on Wander, Fixed Normalise Null Error"
#
# Changes:
#  1. Insert a new sub-bullet "Editor Weight" (ilvl=2, green 00B050) right
#     before the "Arrive" bullet.
#  2. Insert a new sub-bullet "Min speed" (ilvl=2, green 00B050) right
#     before the "Post week 3 tasks to forum with explanations" bullet.
#  3. Mark the existing "Circle moves behind enemy" bullet green (00B050).
#  4. After "Debug lines are outside of circle", add two new sub-bullets:
#       "Circle distance changes over time" (no special colour)
#       "Maths Magnitude 0 check" (green 00B050)

$d = $word.ActiveDocument

# wdColor value for RGB 00B050 (stored as 0x00BBGGRR -> 0x0050B000)
$green = 5287936

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $i
        }
    }
    return -1
}

# --- 4. Two new bullets after "Debug lines are outside of circle" ---
$idx = Find-ParagraphByText $d "Debug lines are outside of circle"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter() | Out-Null

$idx2 = $idx + 1
$p2 = $d.Paragraphs.Item($idx2)
$p2.Range.Text = "Circle distance changes over time"

$p2.Range.InsertParagraphAfter() | Out-Null
$idx3 = $idx2 + 1
$p3 = $d.Paragraphs.Item($idx3)
$p3.Range.Text = "Maths Magnitude 0 check"
$p3.Range.Font.Color = $green

# --- 3. Colour the existing "Circle moves behind enemy" bullet ---
$idxCircle = Find-ParagraphByText $d "Circle moves behind enemy"
$pCircle = $d.Paragraphs.Item($idxCircle)
$pCircle.Range.Font.Color = $green

# --- 2. New bullet "Min speed" before "Post week 3 tasks..." ---
$idxPost = Find-ParagraphByText $d "Post week 3 tasks to forum with explanations"
$pPost = $d.Paragraphs.Item($idxPost)
$pPost.Range.InsertParagraphBefore() | Out-Null
$idxMinSpeed = Find-ParagraphByText $d "Post week 3 tasks to forum with explanations"
$idxMinSpeed = $idxMinSpeed - 1
$pMinSpeed = $d.Paragraphs.Item($idxMinSpeed)
$pMinSpeed.Range.Text = "Min speed"
$pMinSpeed.Range.ListFormat.ListLevelNumber = 3
$pMinSpeed.Range.Font.Color = $green

# --- 1. New bullet "Editor Weight" before "Arrive" ---
$idxArrive = Find-ParagraphByText $d "Arrive"
$pArrive = $d.Paragraphs.Item($idxArrive)
$pArrive.Range.InsertParagraphBefore() | Out-Null
$idxEditorWeight = Find-ParagraphByText $d "Arrive"
$idxEditorWeight = $idxEditorWeight - 1
$pEditorWeight = $d.Paragraphs.Item($idxEditorWeight)
$pEditorWeight.Range.Text = "Editor Weight"
$pEditorWeight.Range.ListFormat.ListLevelNumber = 3
$pEditorWeight.Range.Font.Color = $green

Write-Output "edits applied"
